# MarketBeat-style analyst rank report update:
#  - insert 3 new date columns (Jun_27, Jun_26, Jun_26) before the existing
#    Jun_17/Jun_15/Jun_13/Jun_10 columns
#  - record this week's rating action for Zacks Investment Research
#    (highlighted) in the new Jun_27 column
#  - add a new analyst group: Benchmark / Evercore ISI rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for 3 new date columns in front of the current first data
# column (B). The existing B:E (Jun_17, Jun_15, Jun_13, Jun_10) shift right
# to E:H.
$ws.Columns("B:D").Insert()

# New header row values for the freshly inserted date columns.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# The inserted columns start out as plain "UN" (unchanged) placeholders,
# same as every other data cell in the sheet.
$ws.Range("B2").Value = "UN"
$ws.Range("C2").Value = "UN"
$ws.Range("D2").Value = "UN"

$ws.Range("B3").Value = "UN"
$ws.Range("C3").Value = "UN"
$ws.Range("D3").Value = "UN"

$ws.Range("B4").Value = "UN"
$ws.Range("C4").Value = "UN"
$ws.Range("D4").Value = "UN"

$ws.Range("C5").Value = "UN"
$ws.Range("D5").Value = "UN"

$ws.Range("B6").Value = "UN"
$ws.Range("C6").Value = "UN"
$ws.Range("D6").Value = "UN"

$ws.Range("B7").Value = "UN"
$ws.Range("C7").Value = "UN"
$ws.Range("D7").Value = "UN"

$ws.Range("B8").Value = "UN"
$ws.Range("C8").Value = "UN"
$ws.Range("D8").Value = "UN"

$ws.Range("B9").Value = "UN"
$ws.Range("C9").Value = "UN"
$ws.Range("D9").Value = "UN"

$ws.Range("B10").Value = "UN"
$ws.Range("C10").Value = "UN"
$ws.Range("D10").Value = "UN"

$ws.Range("B11").Value = "UN"
$ws.Range("C11").Value = "UN"
$ws.Range("D11").Value = "UN"

$ws.Range("B12").Value = "UN"
$ws.Range("C12").Value = "UN"
$ws.Range("D12").Value = "UN"

$ws.Range("B13").Value = "UN"
$ws.Range("C13").Value = "UN"
$ws.Range("D13").Value = "UN"

$ws.Range("B14").Value = "UN"
$ws.Range("C14").Value = "UN"
$ws.Range("D14").Value = "UN"

$ws.Range("B15").Value = "UN"
$ws.Range("C15").Value = "UN"
$ws.Range("D15").Value = "UN"

$ws.Range("B16").Value = "UN"
$ws.Range("C16").Value = "UN"
$ws.Range("D16").Value = "UN"

$ws.Range("B17").Value = "UN"
$ws.Range("C17").Value = "UN"
$ws.Range("D17").Value = "UN"

$ws.Range("B18").Value = "UN"
$ws.Range("C18").Value = "UN"
$ws.Range("D18").Value = "UN"

$ws.Range("B19").Value = "UN"
$ws.Range("C19").Value = "UN"
$ws.Range("D19").Value = "UN"

$ws.Range("B20").Value = "UN"
$ws.Range("C20").Value = "UN"
$ws.Range("D20").Value = "UN"

$ws.Range("B21").Value = "UN"
$ws.Range("C21").Value = "UN"
$ws.Range("D21").Value = "UN"

$ws.Range("B22").Value = "UN"
$ws.Range("C22").Value = "UN"
$ws.Range("D22").Value = "UN"

$ws.Range("B23").Value = "UN"
$ws.Range("C23").Value = "UN"
$ws.Range("D23").Value = "UN"

$ws.Range("B24").Value = "UN"
$ws.Range("C24").Value = "UN"
$ws.Range("D24").Value = "UN"

$ws.Range("B25").Value = "UN"
$ws.Range("C25").Value = "UN"
$ws.Range("D25").Value = "UN"

$ws.Range("B26").Value = "UN"
$ws.Range("C26").Value = "UN"
$ws.Range("D26").Value = "UN"

$ws.Range("B27").Value = "UN"
$ws.Range("C27").Value = "UN"
$ws.Range("D27").Value = "UN"

# This week's downgrade for Zacks Investment Research, highlighted like the
# other "new this week" rating-change cells on the sheet.
$ws.Range("B5").Value = "6/27/2018,Downgrades,Buy -> Hold,"
$ws.Range("B5").Interior.Pattern = -4142
$ws.Range("B5").Interior.ColorIndex = 45

# New analyst coverage group appended at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
